$wb = $excel.ActiveWorkbook

# --- Update the "cht13_f_s" sheet (sheet11.xml) with new data values ---
$ws = $wb.Worksheets.Item("cht13_f_s")

$aVals = @(1,1.0250855188141299,1.0250855188141299,1.0490307867730899,1.0752565564424099,1.1094640820980599,1.1436716077536999,1.17331812998859,1.2006841505131101,1.2497149372862,1.22576966932725,1.2884834663625899,1.32155074116305,1.3534777651083201,1.3842645381984,1.41619156214367,1.45039908779931,1.49258836944127,1.5302166476624799,1.56784492588369,1.6088939566704601,1.6510832383124201,1.69783352337514,1.7400228050171,1.77651083238312,1.8198403648802699,1.8677309007981699,1.9144811858608799,1.95496009122006,2)
$bVals = @(0.71054263565891396,0.74775193798449602,0.74775193798449602,0.78248062015503805,0.81596899224806196,0.85813953488371997,0.90031007751937897,0.93379844961240299,0.96480620155038699,1.01751937984496,0.99271317829457295,1.0565891472868201,1.08883720930232,1.1186046511627901,1.1458914728682099,1.17317829457364,1.2017054263565801,1.23395348837209,1.26248062015503,1.2872868217054201,1.31395348837209,1.34062015503875,1.36666666666666,1.38775193798449,1.40511627906976,1.42372093023255,1.44232558139534,1.4572093023255801,1.4689922480620099,1.4807751937984399)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

# apply the "0.000" number format to all populated data cells (A2:B31)
$dataRange = $ws.Range("A2:B31")
$dataRange.NumberFormat = "0.000"

# columns A:B sized to "best fit" their (now wider) numeric content
# (mirrors the bestFit column width already used on the sibling "cht12_f_s" sheet)
$ws.Columns.Item(1).ColumnWidth = 4.65
$ws.Columns.Item(2).ColumnWidth = 4.65

# update selection / active cell on this sheet and make it the active (selected) tab
$ws.Activate()
$ws.Range("A2").Select()

# --- The previously-active sheet "cht12_f_s" (sheet9.xml) is no longer the selected tab ---
$wsPrev = $wb.Worksheets.Item("cht12_f_s")
$wsPrev.Range("C3").Select()

# re-activate cht13_f_s so it ends up as the active/selected sheet and the selection sticks
$ws.Activate()
$ws.Range("A2").Select()

# --- Update workbook window view settings ---
$excel.Windows.Item(1).WindowState = -4143
$excel.Width = 21600
$excel.Height = 11295
$excel.Left = 390
$excel.Top = 390
